# The "Artfynd" (species observation) sheet had its records re-keyed: the
# seven data rows previously on lines 9, 11, 13, 15, 16, 17 and 18 are
# reshuffled so that each row position now holds a different underlying
# observation (new Id in column A, with its species/taxon columns and
# coordinates carried along). No other rows or columns on the sheet change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New content for each affected row, keyed by its row number.
# Columns: A (Id), B (Taxonsorteringsordning), D (Rödlistade), E (TaxonId),
# F (Artnamn), G (Vetenskapligt namn), H (Auktor), Q (Ost), R (Nord).
$rows = @(
    @{ Row = 9;  A = 111670593; B = 78578;  D = "NT"; E = 6458;   F = "Lunglav";   G = "Lobaria pulmonaria"; H = "(L.) Hoffm."; Q = 558040.5475534229; R = 7067901.063021242 },
    @{ Row = 11; A = 111670599; B = 96348;  D = "VU"; E = 220787; F = "Knärot";    G = "Goodyera repens";    H = "(L.) R. Br."; Q = 558031.5226908802; R = 7067909.315233406 },
    @{ Row = 13; A = 111671384; B = 96348;  D = "VU"; E = 220787; F = "Knärot";    G = "Goodyera repens";    H = "(L.) R. Br."; Q = 557798.0632258818; R = 7068181.046264404 },
    @{ Row = 15; A = 111670588; B = 96348;  D = "VU"; E = 220787; F = "Knärot";    G = "Goodyera repens";    H = "(L.) R. Br."; Q = 558039.6361001397; R = 7067902.375451046 },
    @{ Row = 16; A = 111670575; B = 96346;  D = "NT"; E = 620;    F = "Skogsfru";  G = "Epipogium aphyllum"; H = "Sw.";         Q = 558082.6649719321; R = 7067974.943554637 },
    @{ Row = 17; A = 111671406; B = 78578;  D = "NT"; E = 6458;   F = "Lunglav";   G = "Lobaria pulmonaria"; H = "(L.) Hoffm."; Q = 557823.3030943703; R = 7068159.357501161 },
    @{ Row = 18; A = 111671345; B = 96348;  D = "VU"; E = 220787; F = "Knärot";    G = "Goodyera repens";    H = "(L.) R. Br."; Q = 557812.5300353739; R = 7068166.248475613 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
    $ws.Range("Q$n").Value = $r.Q
    $ws.Range("R$n").Value = $r.R
}

# Row 15 now holds a "Knärot" (Goodyera repens) record, which (like the
# other Knärot rows) carries a "Kön" value in column L; row 17 now holds a
# "Lunglav" (Lobaria pulmonaria) record, which never has a column-L entry.
$ws.Range("L15").Value = ""
$ws.Range("L17").ClearContents()
